$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fecha 44382 -> 44281 ; Volumen 160 -> 120 ; Precio min 7000 -> 5500 ; Precio max 8000 -> 6000 ; Precio prom 7438 -> 5750 ; Precio $/Kg 124 -> 96
$ws.Range("D2").Value = 44281
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 5500
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 5750
$ws.Range("P2").Value = 96

# Row 3: Fecha 44281 -> 44421 ; Volumen 120 -> 100 ; Precio min 5500 -> 8000 ; Precio max 6000 -> 9000 ; Precio prom 5750 -> 8500 ; Precio $/Kg 96 -> 142
$ws.Range("D3").Value = 44421
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("P3").Value = 142

# Row 6: Fecha 44421 -> 44400 ; Volumen 100 -> 120 ; Precio min 8000 -> 9000 ; Precio max 9000 -> 10000 ; Precio prom 8500 -> 9500 ; Precio $/Kg 142 -> 158
$ws.Range("D6").Value = 44400
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 9500
$ws.Range("P6").Value = 158

# Row 7: Fecha 44400 -> 44494 ; Volumen 120 (unchanged) ; Precio min 9000 -> 5000 ; Precio max 10000 -> 6000 ; Precio prom 9500 -> 5500 ; Precio $/Kg 158 -> 92
$ws.Range("D7").Value = 44494
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 5500
$ws.Range("P7").Value = 92

# Row 8: Fecha 44494 -> 44382 ; Volumen 120 -> 160 ; Precio min 5000 -> 7000 ; Precio max 6000 -> 8000 ; Precio prom 5500 -> 7438 ; Precio $/Kg 92 -> 124
$ws.Range("D8").Value = 44382
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 7438
$ws.Range("P8").Value = 124
